# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# Cell B2 ("startup" sheet) holds the Cypher query used to populate the
# "CasesTab" output. This updates the query text:
#   1. Removes the blank line that followed the first MATCH clause
#      (...(demo:demographic) / MATCH (c)<--(diag:diagnosis)... now on
#      consecutive lines instead of separated by an empty line).
#   2. Drops the trailing `co.cohort_description` / `Cohort` column from
#      the RETURN clause, so the query now ends after `Response to Treatment`.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`r`nMATCH (c)<--(diag:diagnosis)`r`n MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`r`n`tWHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in  ['T3N1M0', 'Not Applicable']  OPTIONAL MATCH (samp:sample)-->(c)`r`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`r`nWITH DISTINCT c, s, demo, diag, co`r`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`r`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`r`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`r`n        coalesce(demo.breed, '') AS Breed ,`r`n        coalesce(diag.disease_term, '') AS Diagnosis ,`r`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`r`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`r`n        coalesce(demo.sex, '') AS Sex ,`r`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`r`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`r`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# The row auto-sizes to the (now shorter) wrapped text.
$ws.Rows(2).RowHeight = 290

# Leave the same cell selected/active, matching the saved view state.
$ws.Range("B2").Select()
